# Update "想去人数" (interested-people count) figures in the 展览 and 全部类型
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 312
$wsExhibit.Range("F4").Value = 8162
$wsExhibit.Range("F5").Value = 5944
$wsExhibit.Range("F9").Value = 71
$wsExhibit.Range("F11").Value = 671

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 312
$wsAll.Range("F4").Value = 8162
$wsAll.Range("F5").Value = 5944
$wsAll.Range("F9").Value = 71
$wsAll.Range("F15").Value = 671
